# Fruta / hortaliza, semanal
# Insert the new weekly record at row 194 (shifting the existing rows 194-239
# down to 195-240) and populate it with the new week's data for
# Hortaliza / Ajo / Chino / Primera at Terminal La Palmera de La Serena.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 194; everything below shifts down one row
# (old row 194 -> 195, ..., old row 239 -> 240), growing the sheet to R240.
$ws.Rows.Item(194).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(194, 1).Value  = 8
$ws.Cells.Item(194, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(194, 3).Value  = "Coquimbo"
$ws.Cells.Item(194, 4).Value  = 44641
$ws.Cells.Item(194, 5).Value  = 4
$ws.Cells.Item(194, 6).Value  = 100112003
$ws.Cells.Item(194, 7).Value  = "Ajo"
$ws.Cells.Item(194, 8).Value  = "Chino"
$ws.Cells.Item(194, 9).Value  = "Primera"
$ws.Cells.Item(194, 10).Value = 540
$ws.Cells.Item(194, 11).Value = 19000
$ws.Cells.Item(194, 12).Value = 20000
$ws.Cells.Item(194, 13).Value = 19500
$ws.Cells.Item(194, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(194, 15).Value = "China"
$ws.Cells.Item(194, 16).Value = 1950
$ws.Cells.Item(194, 17).Value = 10
$ws.Cells.Item(194, 18).Value = "Hortaliza"
